$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (Sheet1 -> playtest1)
$ws.Name = "playtest1"

# Widen column C for the new long dialogue-text column
$ws.Columns("C:C").ColumnWidth = 41.5

$ws.Range("A17").Value = "ENEMY"
$ws.Range("B17").Value = "EMPTY"
$ws.Range("A18").Value = "ENEMY"
$ws.Range("B18").Value = "Changeling"
$ws.Range("A19").Value = "END_SCENE"
$ws.Range("A20").Value = "START_SCENE"
$ws.Range("B20").Value = "CUTSCENE"
$ws.Range("A21").Value = "DIALOGUE"
$ws.Range("B21").Value = "Clarke"
$ws.Range("C21").Value = "Next Wave incoming!"
$ws.Range("D21").Value = "_"
$ws.Range("E21").Value = "END_DIALOGUE"
$ws.Range("A22").Value = "END_SCENE"
$ws.Range("A23").Value = "START_SCENE"
$ws.Range("B23").Value = "BATTLE"
$ws.Range("A24").Value = "MUSIC"
$ws.Range("B24").Value = "_"
$ws.Range("A25").Value = "ENEMY"
$ws.Range("B25").Value = "Changeling"
$ws.Range("A26").Value = "ENEMY"
$ws.Range("B26").Value = "EMPTY"
$ws.Range("A27").Value = "ENEMY"
$ws.Range("B27").Value = "Tanuki"
$ws.Range("A28").Value = "END_SCENE"
$ws.Range("A29").Value = "START_SCENE"
$ws.Range("B29").Value = "CUTSCENE"
$ws.Range("A30").Value = "DIALOGUE"
$ws.Range("B30").Value = "Clarke"
$ws.Range("C30").Value = "Notha wave comin yo way"
$ws.Range("D30").Value = "_"
$ws.Range("E30").Value = "END_DIALOGUE"
$ws.Range("A31").Value = "END_SCENE"
$ws.Range("A32").Value = "START_SCENE"
$ws.Range("B32").Value = "BATTLE"
$ws.Range("A33").Value = "MUSIC"
$ws.Range("B33").Value = "_"
$ws.Range("A34").Value = "ENEMY"
$ws.Range("B34").Value = "Changeling"
$ws.Range("A35").Value = "ENEMY"
$ws.Range("B35").Value = "Tanuki2"
$ws.Range("A36").Value = "ENEMY"
$ws.Range("B36").Value = "Changeling"
$ws.Range("A37").Value = "END_SCENE"
$ws.Range("A38").Value = "START_SCENE"
$ws.Range("B38").Value = "BATTLE"
$ws.Range("A39").Value = "MUSIC"
$ws.Range("B39").Value = "_"
$ws.Range("A40").Value = "ENEMY"
$ws.Range("B40").Value = "EMPTY"
$ws.Range("A41").Value = "ENEMY"
$ws.Range("B41").Value = "Ijiraq"
$ws.Range("A42").Value = "END_SCENE"
$ws.Range("A43").Value = "START_SCENE"
$ws.Range("B43").Value = "BATTLE"
$ws.Range("A44").Value = "MUSIC"
$ws.Range("B44").Value = "_"
$ws.Range("A45").Value = "ENEMY"
$ws.Range("B45").Value = "Ijiraq2"
$ws.Range("A46").Value = "ENEMY"
$ws.Range("B46").Value = "EMPTY"
$ws.Range("A47").Value = "ENEMY"
$ws.Range("B47").Value = "Tanuki"
$ws.Range("A48").Value = "END_SCENE"
$ws.Range("A49").Value = "START_SCENE"
$ws.Range("B49").Value = "BATTLE"
$ws.Range("A50").Value = "MUSIC"
$ws.Range("B50").Value = "_"
$ws.Range("A51").Value = "ENEMY"
$ws.Range("B51").Value = "Tanuki"
$ws.Range("A52").Value = "ENEMY"
$ws.Range("B52").Value = "Ijiraq2"
$ws.Range("A53").Value = "ENEMY"
$ws.Range("B53").Value = "Changeling"
$ws.Range("A54").Value = "END_SCENE"
$ws.Range("A55").Value = "START_SCENE"
$ws.Range("B55").Value = "BATTLE"
$ws.Range("A56").Value = "MUSIC"
$ws.Range("B56").Value = "bgm_battle_b1"
$ws.Range("A57").Value = "ENEMY"
$ws.Range("B57").Value = "EMPTY"
$ws.Range("A58").Value = "ENEMY"
$ws.Range("B58").Value = "Doppelganger (BLUE)"
$ws.Range("A59").Value = "END_SCENE"
$ws.Range("A60").Value = "END_GAME"

# Restore the on-screen selection to match the authored state
$ws.Range("B30").Select()
